$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 49346
$ws.Range("J87").Value = 49346
$ws.Range("L87").Value = 49346
$ws.Range("N87").Value = -51842
$ws.Range("H90").Value = 49346
$ws.Range("J90").Value = 49346
$ws.Range("L90").Value = 148038
$ws.Range("N90").Value = -160518
$ws.Range("H114").Value = 45710
$ws.Range("J114").Value = 45710
$ws.Range("L114").Value = 45710
$ws.Range("N114").Value = -54388

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 50758.668
$ws.Range("J80").Value = 50758.668
$ws.Range("L80").Value = 50758.668
$ws.Range("N80").Value = -52754.668
$ws.Range("H83").Value = 50758.668
$ws.Range("J83").Value = 50758.668
$ws.Range("L83").Value = 152276.004
$ws.Range("N83").Value = -162260.004
$ws.Range("H111").Value = 40189.2
$ws.Range("J111").Value = 40189.2
$ws.Range("L111").Value = 40189.2
$ws.Range("N111").Value = -48369.2
$ws.Range("H113").Value = 46390
$ws.Range("J113").Value = 46390
$ws.Range("L113").Value = 46390
$ws.Range("N113").Value = -55068
$ws.Range("H114").Value = 41376.668
$ws.Range("J114").Value = 41376.668
$ws.Range("L114").Value = 41376.668
$ws.Range("N114").Value = -50054.668
$ws.Range("H119").Value = 34934.332
$ws.Range("J119").Value = 34934.332
$ws.Range("L119").Value = 34934.332
$ws.Range("N119").Value = -44610.332
$ws.Range("H121").Value = 39996
$ws.Range("J121").Value = 39996
$ws.Range("L121").Value = 39996
$ws.Range("N121").Value = -43490

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 45680
$ws.Range("J108").Value = 45680
$ws.Range("L108").Value = 45680
$ws.Range("N108").Value = -53360
$ws.Range("H110").Value = 48586
$ws.Range("J110").Value = 48586
$ws.Range("L110").Value = 48586
$ws.Range("N110").Value = -56766
$ws.Range("H111").Value = 47702
$ws.Range("J111").Value = 47702
$ws.Range("L111").Value = 47702
$ws.Range("N111").Value = -55882
$ws.Range("H112").Value = 45921.332
$ws.Range("J112").Value = 45921.332
$ws.Range("L112").Value = 45921.332
$ws.Range("N112").Value = -48875.332
$ws.Range("H116").Value = 43490
$ws.Range("J116").Value = 43490
$ws.Range("L116").Value = 43490
$ws.Range("N116").Value = -52668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 4294
$ws.Range("I36").Value = 4492.5
$ws.Range("J36").Value = 3500
$ws.Range("K36").Value = 4492.5
$ws.Range("L36").Value = 3500
$ws.Range("M36").Value = -4104.5
$ws.Range("N36").Value = -4276
$ws.Range("H40").Value = 4294
$ws.Range("I40").Value = 4492.5
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 4492.5
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -4332.5
$ws.Range("N40").Value = -3820
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H82").Value = 37562
$ws.Range("J82").Value = 37562
$ws.Range("L82").Value = 37562
$ws.Range("N82").Value = -38284
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H85").Value = 37562
$ws.Range("J85").Value = 37562
$ws.Range("L85").Value = 37562
$ws.Range("N85").Value = -40058
$ws.Range("H88").Value = 27717.285
$ws.Range("J88").Value = 27717.285
$ws.Range("L88").Value = 27717.285
$ws.Range("N88").Value = -28529.285
$ws.Range("H91").Value = 27717.285
$ws.Range("J91").Value = 27717.285
$ws.Range("L91").Value = 27717.285
$ws.Range("N91").Value = -30525.285
$ws.Range("H110").Value = 42199.5
$ws.Range("J110").Value = 42199.5
$ws.Range("L110").Value = 42199.5
$ws.Range("N110").Value = -50379.5
$ws.Range("H111").Value = 47702
$ws.Range("J111").Value = 47702
$ws.Range("L111").Value = 47702
$ws.Range("N111").Value = -55882
$ws.Range("H112").Value = 29354.2
$ws.Range("J112").Value = 29354.2
$ws.Range("L112").Value = 29354.2
$ws.Range("N112").Value = -32308.2
$ws.Range("H119").Value = 49380
$ws.Range("J119").Value = 49380
$ws.Range("L119").Value = 49380
$ws.Range("N119").Value = -59056

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 21163.924
$ws.Range("J74").Value = 21163.924
$ws.Range("L74").Value = 21163.924
$ws.Range("N74").Value = -23035.924
$ws.Range("H77").Value = 21163.924
$ws.Range("J77").Value = 21163.924
$ws.Range("L77").Value = 63491.772
$ws.Range("N77").Value = -72851.772
$ws.Range("H114").Value = 40330.332
$ws.Range("J114").Value = 40330.332
$ws.Range("L114").Value = 40330.332
$ws.Range("N114").Value = -49008.332
$ws.Range("H116").Value = 40348.4
$ws.Range("J116").Value = 40348.4
$ws.Range("L116").Value = 40348.4
$ws.Range("N116").Value = -49526.4
$ws.Range("H119").Value = 47753
$ws.Range("J119").Value = 47753
$ws.Range("L119").Value = 47753
$ws.Range("N119").Value = -57429
$ws.Range("H123").Value = 15104.333
$ws.Range("J123").Value = 15104.333
$ws.Range("L123").Value = 15104.333
$ws.Range("N123").Value = -20004.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 30813.625
$ws.Range("I99").Value = 20254.5
$ws.Range("J99").Value = 34333.332
$ws.Range("K99").Value = 20254.5
$ws.Range("L99").Value = 34333.332
$ws.Range("M99").Value = -17259.5
$ws.Range("N99").Value = -40323.332
$ws.Range("H102").Value = 40030.5
$ws.Range("J102").Value = 40030.5
$ws.Range("L102").Value = 40030.5
$ws.Range("N102").Value = -46520.5
$ws.Range("H108").Value = 48626
$ws.Range("J108").Value = 48626
$ws.Range("L108").Value = 48626
$ws.Range("N108").Value = -56306
$ws.Range("H109").Value = 35281
$ws.Range("J109").Value = 35281
$ws.Range("L109").Value = 35281
$ws.Range("N109").Value = -38055
$ws.Range("H112").Value = 35839.5
$ws.Range("J112").Value = 35839.5
$ws.Range("L112").Value = 35839.5
$ws.Range("N112").Value = -38793.5
$ws.Range("H114").Value = 39394
$ws.Range("J114").Value = 39394
$ws.Range("L114").Value = 39394
$ws.Range("N114").Value = -48072
$ws.Range("H116").Value = 45664
$ws.Range("J116").Value = 45664
$ws.Range("L116").Value = 45664
$ws.Range("N116").Value = -54842
$ws.Range("H117").Value = 35325.332
$ws.Range("J117").Value = 35325.332
$ws.Range("L117").Value = 35325.332
$ws.Range("N117").Value = -44503.332
$ws.Range("H118").Value = 38362
$ws.Range("J118").Value = 38362
$ws.Range("L118").Value = 38362
$ws.Range("N118").Value = -41676
$ws.Range("H119").Value = 47404
$ws.Range("J119").Value = 47404
$ws.Range("L119").Value = 47404
$ws.Range("N119").Value = -57080
$ws.Range("H120").Value = 46391.6
$ws.Range("J120").Value = 46391.6
$ws.Range("L120").Value = 46391.6
$ws.Range("N120").Value = -56067.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 28200
$ws.Range("J87").Value = 28200
$ws.Range("L87").Value = 28200
$ws.Range("N87").Value = -30696
$ws.Range("H90").Value = 28200
$ws.Range("J90").Value = 28200
$ws.Range("L90").Value = 84600
$ws.Range("N90").Value = -97080
$ws.Range("H92").Value = 29783.334
$ws.Range("J92").Value = 29783.334
$ws.Range("L92").Value = 29783.334
$ws.Range("N92").Value = -34775.334
$ws.Range("H93").Value = 35914.285
$ws.Range("J93").Value = 35914.285
$ws.Range("L93").Value = 35914.285
$ws.Range("N93").Value = -40906.285
$ws.Range("H99").Value = 38143
$ws.Range("J99").Value = 37777.555
$ws.Range("L99").Value = 37777.555
$ws.Range("N99").Value = -43767.555
$ws.Range("H102").Value = 41337
$ws.Range("J102").Value = 41337
$ws.Range("L102").Value = 41337
$ws.Range("N102").Value = -47827
$ws.Range("H106").Value = 34912
$ws.Range("J106").Value = 34912
$ws.Range("L106").Value = 34912
$ws.Range("N106").Value = -37436
$ws.Range("H108").Value = 30155
$ws.Range("J108").Value = 30155
$ws.Range("L108").Value = 30155
$ws.Range("N108").Value = -37835
$ws.Range("H109").Value = 29892
$ws.Range("J109").Value = 29892
$ws.Range("L109").Value = 29892
$ws.Range("N109").Value = -32666
$ws.Range("H116").Value = 49672
$ws.Range("J116").Value = 49672
$ws.Range("L116").Value = 49672
$ws.Range("N116").Value = -58850
$ws.Range("H117").Value = 47290
$ws.Range("J117").Value = 47290
$ws.Range("L117").Value = 47290
$ws.Range("N117").Value = -56468
$ws.Range("H121").Value = 36954
$ws.Range("J121").Value = 36954
$ws.Range("L121").Value = 36954
$ws.Range("N121").Value = -40448
